$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the two mis-typed names that were corrected in this revision:
#  - D14: "Eszenzi Cilla" -> "Eszenyi Cilla"
#  - E21: "Kvács Éva"     -> "Kovács Éva"
$ws.Range("D14").Value = "Pelák Olgi`nEszenyi Cilla,Gál Edit"
$ws.Range("E21").Value = "Kovács István`nKovács Éva"

# Restore the view state saved with the workbook (scroll position + selection).
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("I16").Select()
